# Weekly update: a new Jengibre price record (2023-07-28) is inserted at
# row 49, pushing the previously-existing rows 49:60 down to rows 50:61.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 49, shifting existing rows 49:60 down to 50:61
$ws.Rows("49:49").Insert()

# Populate the new row 49 with the new data record
$ws.Range("A49").Value = 11
$ws.Range("B49").Value = "Vega Monumental Concepción"
$ws.Range("C49").Value = "Bíobío"
$ws.Range("D49").Value = 45135
$ws.Range("D49").NumberFormat = $ws.Range("D50").NumberFormat
$ws.Range("E49").Value = 8
$ws.Range("F49").Value = 100114007
$ws.Range("G49").Value = "Jengibre"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 30
$ws.Range("K49").Value = 18000
$ws.Range("L49").Value = 18000
$ws.Range("M49").Value = 18000
$ws.Range("N49").Value = '$/caja 13 kilos'
$ws.Range("O49").Value = "Perú"
$ws.Range("P49").Value = 1385
$ws.Range("Q49").Value = 13
$ws.Range("R49").Value = "Hortaliza"
